# Apply the cryptos list update (Thu Mar 21 14:37:59 UTC 2024 GitHub Actions run).
# Values in columns D (Price) and E (Volume 1h) are stored as plain text in this
# workbook (inline strings), including numeric-looking values like "1.00" or
# "0.635" -- so we force text type via NumberFormat="@" before assigning, then
# restore the default "Normal" style so no stray number formatting is left behind.
function Set-CellText {
    param($ws, $addr, $val)
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-CellText $ws 'D2' '67.056.88'
Set-CellText $ws 'E2' '  +4.60%  '
Set-CellText $ws 'D3' '3.543.10'
Set-CellText $ws 'E3' '  +5.48%  '
Set-CellText $ws 'D4' '1.00'
Set-CellText $ws 'E4' '  -0.04%  '
Set-CellText $ws 'D5' '189.00'
Set-CellText $ws 'E5' '  +8.11%  '
Set-CellText $ws 'D6' '563.59'
Set-CellText $ws 'E6' '  +6.54%  '
Set-CellText $ws 'D7' '0.624'
Set-CellText $ws 'E7' '  +4.43%  '
Set-CellText $ws 'D8' '3.539.08'
Set-CellText $ws 'E8' '  +5.41%  '
Set-CellText $ws 'E9' '  -0.05%  '
Set-CellText $ws 'D10' '0.635'
Set-CellText $ws 'E10' '  +4.03%  '
Set-CellText $ws 'D11' '0.155'
Set-CellText $ws 'E11' '  +15.14%  '
Set-CellText $ws 'D12' '54.87'
Set-CellText $ws 'E12' '  +2.49%  '
Set-CellText $ws 'D13' '0.0000275'
Set-CellText $ws 'E13' '  +6.66%  '
Set-CellText $ws 'D14' '9.36'
Set-CellText $ws 'E14' '  +2.74%  '
Set-CellText $ws 'D15' '4.107.44'
Set-CellText $ws 'E15' '  +5.54%  '
Set-CellText $ws 'D16' '3.551.37'
Set-CellText $ws 'E16' '  +6.06%  '
Set-CellText $ws 'D17' '18.62'
Set-CellText $ws 'E17' '  +5.71%  '
Set-CellText $ws 'E18' '  +3.32%  '
Set-CellText $ws 'D19' '67.159.43'
Set-CellText $ws 'E19' '  +4.74%  '
Set-CellText $ws 'D20' '12.10'
Set-CellText $ws 'E20' '  +7.74%  '
Set-CellText $ws 'E21' '  +3.37%  '
Set-CellText $ws 'D22' '422.40'
Set-CellText $ws 'E22' '  +12.83%  '
Set-CellText $ws 'D23' '4.13'
Set-CellText $ws 'E23' '  +10.65%  '
Set-CellText $ws 'D24' '85.41'
Set-CellText $ws 'E25' '  +1.92%  '
Set-CellText $ws 'D26' '11.07'
Set-CellText $ws 'E26' '  -4.71%  '
Set-CellText $ws 'D27' '2.92'
Set-CellText $ws 'E27' '  +7.95%  '
Set-CellText $ws 'D28' '12.33'
Set-CellText $ws 'E28' '  +8.67%  '
Set-CellText $ws 'E29' '  -1.09%  '
Set-CellText $ws 'E30' '  +10.18%  '
Set-CellText $ws 'D31' '30.48'
Set-CellText $ws 'E31' '  +5.31%  '
Set-CellText $ws 'D32' '632.22'
Set-CellText $ws 'E32' '  -0.56%  '
Set-CellText $ws 'D33' '6.68'
Set-CellText $ws 'E33' '  +3.03%  '
Set-CellText $ws 'D34' '11.76'
Set-CellText $ws 'E34' '  +4.60%  '
Set-CellText $ws 'E35' '  +5.07%  '
Set-CellText $ws 'D36' '60.38'
Set-CellText $ws 'E36' '  +4.13%  '
Set-CellText $ws 'D37' '0.0₃0836'
Set-CellText $ws 'E37' '  +14.60%  '
Set-CellText $ws 'B38' 'Kaspa'
Set-CellText $ws 'C38' 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-CellText $ws 'D38' '0.149'
Set-CellText $ws 'E38' '  +19.25%  '
Set-CellText $ws 'B39' 'InjectiveProtocol'
Set-CellText $ws 'C39' 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-CellText $ws 'D39' '38.42'
Set-CellText $ws 'E39' '  +5.09%  '
Set-CellText $ws 'E40' '  -0.12%  '
Set-CellText $ws 'D41' '0.388'
Set-CellText $ws 'E41' '  +1.84%  '
Set-CellText $ws 'B42' 'Maker'
Set-CellText $ws 'C42' 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-CellText $ws 'D42' '3.153.37'
Set-CellText $ws 'E42' '  +5.88%  '
Set-CellText $ws 'B43' 'Stacks'
Set-CellText $ws 'C43' 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-CellText $ws 'D43' '3.36'
Set-CellText $ws 'E43' '  +11.12%  '
Set-CellText $ws 'D44' '1.00'
Set-CellText $ws 'E44' '  +0.05%  '
Set-CellText $ws 'D45' '2.64'
Set-CellText $ws 'E45' '  -0.92%  '
Set-CellText $ws 'D46' '2.87'
Set-CellText $ws 'E46' '  +9.56%  '
Set-CellText $ws 'D47' '3.37'
Set-CellText $ws 'E47' '  +10.80%  '
Set-CellText $ws 'D48' '0.0419'
Set-CellText $ws 'E48' '  +5.29%  '
Set-CellText $ws 'E49' '  +1.97%  '
Set-CellText $ws 'D50' '0.132'
Set-CellText $ws 'E50' '  +5.29%  '
Set-CellText $ws 'D51' '8.63'
Set-CellText $ws 'E51' '  +8.13%  '
